# repull data, push all data, mean calculation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F ("dSF") values were refreshed after repulling data.
$ws.Range("F2").Value = 3
$ws.Range("F10").Value = -3
$ws.Range("F15").Value = 2
$ws.Range("F19").Value = -5
$ws.Range("F25").Value = -2
